$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Label" in new column H, matching the style of the neighboring header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Fill in the new "Label" column values for the two blocks of rows (2-11 and 12-21)
$labelValues = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labelValues[$i]
}

for ($i = 0; $i -lt 10; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 8).Value = $labelValues[$i]
}
